# Updates the "想去人数" (F column) counts on the "展览" and "全部类型"
# worksheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F (想去人数)
$updates = @{
    2  = 307
    4  = 87
    5  = 388
    6  = 11412
    7  = 711
    8  = 108
    10 = 84
    12 = 161
    13 = 20
    19 = 1276
    21 = 894
    22 = 109
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
